$d = $word.ActiveDocument

# The site rebuild dropped the trailing "Ver no Jupiter Salvar em pdf Salvar
# em docx" line, the copyright/footer line right after it, and the blank
# paragraph that separated them from the preceding "Requisitos" section
# (the paragraph that reads "LOQ4031: Química Geral I (Requisito fraco)").
#
# Find the first paragraph of that trailing block ("Ver no Jupiter ...") and
# its last paragraph (the "© 2020 ..." copyright line). Everything from the
# start of the blank paragraph right before the block through to the end of
# the copyright paragraph (paragraph mark included) gets removed, so the
# paragraphs disappear rather than merging with their neighbours.

$startMarker = "Ver no Jupiter Salvar em pdf Salvar em docx"
$endMarker = "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$firstPara = $null
$lastPara = $null

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if (($firstPara -eq $null) -and ($text -eq $startMarker)) {
        $firstPara = $p
    }
    if ($text -eq $endMarker) {
        $lastPara = $p
    }
}

if (($firstPara -ne $null) -and ($lastPara -ne $null)) {
    $deleteStart = $firstPara.Previous().Range.Start
    $deleteEnd = $lastPara.Range.End
    $r = $d.Range($deleteStart, $deleteEnd)
    $r.Delete()
}
